$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: DeiT header
$ws.Range("A14").Value = "DeiT"

# Row 15: STL  sub-header
$ws.Range("A15").Value = "STL "

# Row 16: Training Loss (STL)
$ws.Range("A16").Value = "Training Loss"
$vals16 = @(2.38142881267949,2.3402747919684899,2.0297876063146001,1.5991922303249899,1.4425443912807201,1.36741580147492,1.2814129779213299,1.3011209149109599,1.3564058244228301,1.2828656105618701,1.2718531881508,1.2773549619473901,1.1894847154617301,1.2151779124611299,1.24914535252671,1.2079322557700301,1.21191795248734,1.19468530228263,1.1078823381348599,1.20008581249337)
for ($i = 0; $i -lt $vals16.Length; $i++) { $ws.Cells.Item(16, $i + 2).Value = $vals16[$i] }

# Row 17: Test accuracy (STL)
$ws.Range("A17").Value = "Test accuracy"
$vals17 = @(12.0625002365112,14.7625003967285,76.387503173828094,89.412502990722601,92.600002441406204,94.437502746581998,94.175002380371097,94.137502441406198,95.012501953124996,95.350002624511703,95.162502807617102,95.825002746582001,95.687502197265601,95.725002746582007,95.762502624511697,95.812502197265601,95.950002319335894,96.112502197265599,96.162502380371095,96.162502075195306)
for ($i = 0; $i -lt $vals17.Length; $i++) { $ws.Cells.Item(17, $i + 2).Value = $vals17[$i] }

# Row 19: Caltech101 sub-header
$ws.Range("A19").Value = "Caltech101"

# Row 20: Training Loss (Caltech101)
$ws.Range("A20").Value = "Training Loss"
$vals20 = @(4.7773776572683504,4.7209026398866003,4.3020073382750796,3.6252780530763702,3.1442557832469098,2.7128629269807201,2.3086596260899999,2.148361750271,1.96025768570278,1.9582724182502,1.7233908746553499,1.74969106912612,1.64456712422163,1.57949196514876,1.58193930853968,1.5167814778244999,1.5042027740374799,1.4967302990996301,1.4715914739214799,1.43566106842911)
for ($i = 0; $i -lt $vals20.Length; $i++) { $ws.Cells.Item(20, $i + 2).Value = $vals20[$i] }

# Row 21: Test accuracy (Caltech101)
$ws.Range("A21").Value = "Test accuracy"
$vals21 = @(0.61443934984470805,0.61443934984470805,27.726574875792,44.162827838950399,64.746545800415603,69.508449958766107,76.958528263777595,78.955455024121505,78.571430821572505,79.953919300835196,80.952384913571905,83.486945789530495,83.640557706630702,82.872507051388794,84.024579448084594,84.254994405579396,84.024580116096104,83.717361204085805,84.792630314277602,84.331798360094993)
for ($i = 0; $i -lt $vals21.Length; $i++) { $ws.Cells.Item(21, $i + 2).Value = $vals21[$i] }

# Row 23: CIFAR-10 sub-header
$ws.Range("A23").Value = "CIFAR-10"

# Row 24: Training Loss (CIFAR-10)
$ws.Range("A24").Value = "Training Loss"
$vals24 = @(2.2708631918980502,2.10866297422311,1.69960997807673,1.5880840501724101,1.52696154316266,1.5243850996860999,1.5362521422214901,1.46570869668936,1.4676731491700099,1.46148274005987,1.42804607550303,1.4317569870215101,1.40641588828502,1.40195413292982,1.3760134472296699,1.3802266029211101,1.3791513161781499,1.3245098438018399,1.3455636903261501,1.32338791153369)
for ($i = 0; $i -lt $vals24.Length; $i++) { $ws.Cells.Item(24, $i + 2).Value = $vals24[$i] }

# Row 25: Test accuracy (CIFAR-10)
$ws.Range("A25").Value = "Test accuracy"
$vals25 = @(36.550001074218699,50.830001293945301,91.220002490234293,93.160002783203097,93.9600025878906,93.0200026855468,93.440002392578094,94.270002343749994,94.640002246093701,94.540002343750004,95.750002685546804,95.640002587890606,95.980002197265605,96.290002539062499,96.350001953125002,96.440002685546801,96.670002929687499,96.650002685546795,96.620002832031204,96.900002294921805)
for ($i = 0; $i -lt $vals25.Length; $i++) { $ws.Cells.Item(25, $i + 2).Value = $vals25[$i] }

# Update selection to match target state
$ws.Range("K30").Select()